$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 11-13 (the "Resolving-Mac" sending-cluster block), shifting
# the table dimension from A1:T13 down to A1:T10.
$ws.Rows("11:13").Delete()

# Update recalculated TPM-derived numeric values for rows 2-10.
# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1.471704
$ws.Range("H2").Value2 = 4.415112
$ws.Range("I2").Value2 = 0.08657490103749592
$ws.Range("J2").Value2 = 0.0865749010374959
$ws.Range("M2").Value2 = 0.02270466666666667
$ws.Range("N2").Value2 = 0.06811400000000001
$ws.Range("O2").Value2 = 0.002206225855740089
$ws.Range("P2").Value2 = 0.002206225855740089
$ws.Range("Q2").Value2 = 0.033414548752
$ws.Range("R2").Value2 = 0.300730938768
$ws.Range("S2").Value2 = 0.0001910037851270629
$ws.Range("T2").Value2 = 0.0001910037851270629
# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1.471704
$ws.Range("H3").Value2 = 4.415112
$ws.Range("I3").Value2 = 0.08657490103749592
$ws.Range("J3").Value2 = 0.0865749010374959
$ws.Range("O3").Value2 = 0.002281111990432972
$ws.Range("P3").Value2 = 0.002281111990432972
$ws.Range("Q3").Value2 = 0.034548741968
$ws.Range("R3").Value2 = 0.310938677712
$ws.Range("S3").Value2 = 0.0001974870448271799
$ws.Range("T3").Value2 = 0.0001974870448271799
# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1.471704
$ws.Range("H4").Value2 = 4.415112
$ws.Range("I4").Value2 = 0.08657490103749592
$ws.Range("J4").Value2 = 0.0865749010374959
$ws.Range("M4").Value2 = 10.24499966666667
$ws.Range("N4").Value2 = 30.734999
$ws.Range("O4").Value2 = 0.9955126621538269
$ws.Range("P4").Value2 = 0.9955126621538269
$ws.Range("Q4").Value2 = 15.077606989432
$ws.Range("R4").Value2 = 135.698462904888
$ws.Range("S4").Value2 = 0.08618641020754168
$ws.Range("T4").Value2 = 0.08618641020754167
# Row 5
$ws.Range("I5").Value2 = 0.2867067974456365
$ws.Range("J5").Value2 = 0.2867067974456365
$ws.Range("M5").Value2 = 0.02270466666666667
$ws.Range("N5").Value2 = 0.06811400000000001
$ws.Range("O5").Value2 = 0.002206225855740089
$ws.Range("P5").Value2 = 0.002206225855740089
$ws.Range("Q5").Value2 = 0.1106576865346667
$ws.Range("R5").Value2 = 0.9959191788120002
$ws.Range("S5").Value2 = 0.0006325399495409998
$ws.Range("T5").Value2 = 0.0006325399495409998
# Row 6
$ws.Range("I6").Value2 = 0.2867067974456365
$ws.Range("J6").Value2 = 0.2867067974456365
$ws.Range("O6").Value2 = 0.002281111990432972
$ws.Range("P6").Value2 = 0.002281111990432972
$ws.Range("S6").Value2 = 0.000654010313391879
$ws.Range("T6").Value2 = 0.000654010313391879
# Row 7
$ws.Range("I7").Value2 = 0.2867067974456365
$ws.Range("J7").Value2 = 0.2867067974456365
$ws.Range("M7").Value2 = 10.24499966666667
$ws.Range("N7").Value2 = 30.734999
$ws.Range("O7").Value2 = 0.9955126621538269
$ws.Range("P7").Value2 = 0.9955126621538269
$ws.Range("Q7").Value2 = 49.93193594540467
$ws.Range("R7").Value2 = 449.387423508642
$ws.Range("S7").Value2 = 0.2854202471827036
$ws.Range("T7").Value2 = 0.2854202471827036
# Row 8
$ws.Range("G8").Value2 = 10.65370933333333
$ws.Range("H8").Value2 = 31.961128
$ws.Range("I8").Value2 = 0.6267183015168676
$ws.Range("J8").Value2 = 0.6267183015168675
$ws.Range("M8").Value2 = 0.02270466666666667
$ws.Range("N8").Value2 = 0.06811400000000001
$ws.Range("O8").Value2 = 0.002206225855740089
$ws.Range("P8").Value2 = 0.002206225855740089
$ws.Range("Q8").Value2 = 0.2418889191768889
$ws.Range("R8").Value2 = 2.177000272592001
$ws.Range("S8").Value2 = 0.001382682121072026
$ws.Range("T8").Value2 = 0.001382682121072026
# Row 9
$ws.Range("G9").Value2 = 10.65370933333333
$ws.Range("H9").Value2 = 31.961128
$ws.Range("I9").Value2 = 0.6267183015168676
$ws.Range("J9").Value2 = 0.6267183015168675
$ws.Range("O9").Value2 = 0.002281111990432972
$ws.Range("P9").Value2 = 0.002281111990432972
$ws.Range("Q9").Value2 = 0.2500993778364445
$ws.Range("R9").Value2 = 2.250894400528
$ws.Range("S9").Value2 = 0.001429614632213913
$ws.Range("T9").Value2 = 0.001429614632213913
# Row 10
$ws.Range("G10").Value2 = 10.65370933333333
$ws.Range("H10").Value2 = 31.961128
$ws.Range("I10").Value2 = 0.6267183015168676
$ws.Range("J10").Value2 = 0.6267183015168675
$ws.Range("M10").Value2 = 10.24499966666667
$ws.Range("N10").Value2 = 30.734999
$ws.Range("O10").Value2 = 0.9955126621538269
$ws.Range("P10").Value2 = 0.9955126621538269
$ws.Range("Q10").Value2 = 109.1472485687636
$ws.Range("R10").Value2 = 982.3252371188721
$ws.Range("S10").Value2 = 0.6239060047635816
$ws.Range("T10").Value2 = 0.6239060047635815
